$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.0005253333333333333
$ws.Cells.Item(2, 8).Value = 0.001576
$ws.Cells.Item(2, 9).Value = 0.001186555780920845
$ws.Cells.Item(2, 10).Value = 0.001186555780920846
$ws.Cells.Item(2, 13).Value = 0.2057716666666667
$ws.Cells.Item(2, 14).Value = 0.6173149999999999
$ws.Cells.Item(2, 15).Value = 0.01089677771948535
$ws.Cells.Item(2, 16).Value = 0.01089677771948535
$ws.Cells.Item(2, 17).Value = 0.0001080987155555555
$ws.Cells.Item(2, 18).Value = 0.0009728884399999999
$ws.Cells.Item(2, 19).Value = [double]"1.29296345964648E-05"
$ws.Cells.Item(2, 20).Value = [double]"1.29296345964648E-05"
$ws.Cells.Item(3, 7).Value = 0.0005253333333333333
$ws.Cells.Item(3, 8).Value = 0.001576
$ws.Cells.Item(3, 9).Value = 0.001186555780920845
$ws.Cells.Item(3, 10).Value = 0.001186555780920846
$ws.Cells.Item(3, 15).Value = 0.01769706320706529
$ws.Cells.Item(3, 16).Value = 0.01769706320706529
$ws.Cells.Item(3, 17).Value = 0.0001755592204444444
$ws.Cells.Item(3, 18).Value = 0.001580032984
$ws.Cells.Item(3, 19).Value = [double]"2.099855265366491E-05"
$ws.Cells.Item(3, 20).Value = [double]"2.099855265366492E-05"
$ws.Cells.Item(4, 7).Value = 0.0005253333333333333
$ws.Cells.Item(4, 8).Value = 0.001576
$ws.Cells.Item(4, 9).Value = 0.001186555780920845
$ws.Cells.Item(4, 10).Value = 0.001186555780920846
$ws.Cells.Item(4, 13).Value = 0.4895776666666666
$ws.Cells.Item(4, 14).Value = 1.468733
$ws.Cells.Item(4, 15).Value = 0.02592591631545138
$ws.Cells.Item(4, 16).Value = 0.02592591631545138
$ws.Cells.Item(4, 17).Value = 0.0002571914675555555
$ws.Cells.Item(4, 18).Value = 0.002314723207999999
$ws.Cells.Item(4, 19).Value = [double]"3.07625458797689E-05"
$ws.Cells.Item(4, 20).Value = [double]"3.076254587976891E-05"
$ws.Cells.Item(5, 7).Value = 0.0005253333333333333
$ws.Cells.Item(5, 8).Value = 0.001576
$ws.Cells.Item(5, 9).Value = 0.001186555780920845
$ws.Cells.Item(5, 10).Value = 0.001186555780920846
$ws.Cells.Item(5, 13).Value = 17.85418133333333
$ws.Cells.Item(5, 14).Value = 53.562544
$ws.Cells.Item(5, 15).Value = 0.9454802427579979
$ws.Cells.Item(5, 16).Value = 0.945480242757998
$ws.Cells.Item(5, 17).Value = 0.009379396593777777
$ws.Cells.Item(5, 18).Value = 0.084414569344
$ws.Cells.Item(5, 19).Value = 0.001121865047790947
$ws.Cells.Item(5, 20).Value = 0.001121865047790947
$ws.Cells.Item(6, 9).Value = 0.03480764394894197
$ws.Cells.Item(6, 10).Value = 0.03480764394894197
$ws.Cells.Item(6, 13).Value = 0.2057716666666667
$ws.Cells.Item(6, 14).Value = 0.6173149999999999
$ws.Cells.Item(6, 15).Value = 0.01089677771948535
$ws.Cells.Item(6, 16).Value = 0.01089677771948535
$ws.Cells.Item(6, 17).Value = 0.003171078564444445
$ws.Cells.Item(6, 18).Value = 0.02853970708
$ws.Cells.Item(6, 19).Value = 0.0003792911590506097
$ws.Cells.Item(6, 20).Value = 0.0003792911590506097
$ws.Cells.Item(7, 9).Value = 0.03480764394894197
$ws.Cells.Item(7, 10).Value = 0.03480764394894197
$ws.Cells.Item(7, 15).Value = 0.01769706320706529
$ws.Cells.Item(7, 16).Value = 0.01769706320706529
$ws.Cells.Item(7, 19).Value = 0.0006159930750534496
$ws.Cells.Item(7, 20).Value = 0.0006159930750534497
$ws.Cells.Item(8, 9).Value = 0.03480764394894197
$ws.Cells.Item(8, 10).Value = 0.03480764394894197
$ws.Cells.Item(8, 13).Value = 0.4895776666666666
$ws.Cells.Item(8, 14).Value = 1.468733
$ws.Cells.Item(8, 15).Value = 0.02592591631545138
$ws.Cells.Item(8, 16).Value = 0.02592591631545138
$ws.Cells.Item(8, 17).Value = 0.007544718228444445
$ws.Cells.Item(8, 18).Value = 0.06790246405599999
$ws.Cells.Item(8, 19).Value = 0.0009024200641582971
$ws.Cells.Item(8, 20).Value = 0.0009024200641582971
$ws.Cells.Item(9, 9).Value = 0.03480764394894197
$ws.Cells.Item(9, 10).Value = 0.03480764394894197
$ws.Cells.Item(9, 13).Value = 17.85418133333333
$ws.Cells.Item(9, 14).Value = 53.562544
$ws.Cells.Item(9, 15).Value = 0.9454802427579979
$ws.Cells.Item(9, 16).Value = 0.945480242757998
$ws.Cells.Item(9, 17).Value = 0.2751448371342222
$ws.Cells.Item(9, 18).Value = 2.476303534208
$ws.Cells.Item(9, 19).Value = 0.03290993965067961
$ws.Cells.Item(9, 20).Value = 0.03290993965067961
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.09537533333333333
$ws.Cells.Item(10, 8).Value = 0.286126
$ws.Cells.Item(10, 9).Value = 0.2154216112764961
$ws.Cells.Item(10, 10).Value = 0.2154216112764961
$ws.Cells.Item(10, 13).Value = 0.2057716666666667
$ws.Cells.Item(10, 14).Value = 0.6173149999999999
$ws.Cells.Item(10, 15).Value = 0.01089677771948535
$ws.Cells.Item(10, 16).Value = 0.01089677771948535
$ws.Cells.Item(10, 17).Value = 0.01962554129888889
$ws.Cells.Item(10, 18).Value = 0.17662987169
$ws.Cells.Item(10, 19).Value = 0.002347401414053356
$ws.Cells.Item(10, 20).Value = 0.002347401414053356
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.09537533333333333
$ws.Cells.Item(11, 8).Value = 0.286126
$ws.Cells.Item(11, 9).Value = 0.2154216112764961
$ws.Cells.Item(11, 10).Value = 0.2154216112764961
$ws.Cells.Item(11, 15).Value = 0.01769706320706529
$ws.Cells.Item(11, 16).Value = 0.01769706320706529
$ws.Cells.Item(11, 17).Value = 0.03187313293711111
$ws.Cells.Item(11, 18).Value = 0.286858196434
$ws.Cells.Item(11, 19).Value = 0.003812329870928
$ws.Cells.Item(11, 20).Value = 0.003812329870928
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.09537533333333333
$ws.Cells.Item(12, 8).Value = 0.286126
$ws.Cells.Item(12, 9).Value = 0.2154216112764961
$ws.Cells.Item(12, 10).Value = 0.2154216112764961
$ws.Cells.Item(12, 13).Value = 0.4895776666666666
$ws.Cells.Item(12, 14).Value = 1.468733
$ws.Cells.Item(12, 15).Value = 0.02592591631545138
$ws.Cells.Item(12, 16).Value = 0.02592591631545138
$ws.Cells.Item(12, 17).Value = 0.04669363315088888
$ws.Cells.Item(12, 18).Value = 0.4202426983579999
$ws.Cells.Item(12, 19).Value = 0.005585002666494135
$ws.Cells.Item(12, 20).Value = 0.005585002666494135
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.09537533333333333
$ws.Cells.Item(13, 8).Value = 0.286126
$ws.Cells.Item(13, 9).Value = 0.2154216112764961
$ws.Cells.Item(13, 10).Value = 0.2154216112764961
$ws.Cells.Item(13, 13).Value = 17.85418133333333
$ws.Cells.Item(13, 14).Value = 53.562544
$ws.Cells.Item(13, 15).Value = 0.9454802427579979
$ws.Cells.Item(13, 16).Value = 0.945480242757998
$ws.Cells.Item(13, 17).Value = 1.702848496060444
$ws.Cells.Item(13, 18).Value = 15.325636464544
$ws.Cells.Item(13, 19).Value = 0.2036768773250206
$ws.Cells.Item(13, 20).Value = 0.2036768773250206
$ws.Cells.Item(14, 7).Value = 0.3314266666666667
$ws.Cells.Item(14, 8).Value = 0.9942800000000001
$ws.Cells.Item(14, 9).Value = 0.7485841889936411
$ws.Cells.Item(14, 10).Value = 0.7485841889936411
$ws.Cells.Item(14, 13).Value = 0.2057716666666667
$ws.Cells.Item(14, 14).Value = 0.6173149999999999
$ws.Cells.Item(14, 15).Value = 0.01089677771948535
$ws.Cells.Item(14, 16).Value = 0.01089677771948535
$ws.Cells.Item(14, 17).Value = 0.06819821757777778
$ws.Cells.Item(14, 18).Value = 0.6137839582
$ws.Cells.Item(14, 19).Value = 0.008157155511784916
$ws.Cells.Item(14, 20).Value = 0.008157155511784916
$ws.Cells.Item(15, 7).Value = 0.3314266666666667
$ws.Cells.Item(15, 8).Value = 0.9942800000000001
$ws.Cells.Item(15, 9).Value = 0.7485841889936411
$ws.Cells.Item(15, 10).Value = 0.7485841889936411
$ws.Cells.Item(15, 15).Value = 0.01769706320706529
$ws.Cells.Item(15, 16).Value = 0.01769706320706529
$ws.Cells.Item(15, 17).Value = 0.1107582625022222
$ws.Cells.Item(15, 18).Value = 0.9968243625200001
$ws.Cells.Item(15, 19).Value = 0.01324774170843017
$ws.Cells.Item(15, 20).Value = 0.01324774170843018
$ws.Cells.Item(16, 7).Value = 0.3314266666666667
$ws.Cells.Item(16, 8).Value = 0.9942800000000001
$ws.Cells.Item(16, 9).Value = 0.7485841889936411
$ws.Cells.Item(16, 10).Value = 0.7485841889936411
$ws.Cells.Item(16, 13).Value = 0.4895776666666666
$ws.Cells.Item(16, 14).Value = 1.468733
$ws.Cells.Item(16, 15).Value = 0.02592591631545138
$ws.Cells.Item(16, 16).Value = 0.02592591631545138
$ws.Cells.Item(16, 17).Value = 0.1622590941377778
$ws.Cells.Item(16, 18).Value = 1.46033184724
$ws.Cells.Item(16, 19).Value = 0.01940773103891918
$ws.Cells.Item(16, 20).Value = 0.01940773103891918
$ws.Cells.Item(17, 7).Value = 0.3314266666666667
$ws.Cells.Item(17, 8).Value = 0.9942800000000001
$ws.Cells.Item(17, 9).Value = 0.7485841889936411
$ws.Cells.Item(17, 10).Value = 0.7485841889936411
$ws.Cells.Item(17, 13).Value = 17.85418133333333
$ws.Cells.Item(17, 14).Value = 53.562544
$ws.Cells.Item(17, 15).Value = 0.9454802427579979
$ws.Cells.Item(17, 16).Value = 0.945480242757998
$ws.Cells.Item(17, 17).Value = 5.91735180536889
$ws.Cells.Item(17, 18).Value = 53.25616624832001
$ws.Cells.Item(17, 19).Value = 0.7077715607345068
$ws.Cells.Item(17, 20).Value = 0.7077715607345069
